$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 84

$ws.Cells.Item($row, 1).Value = "EKSN87"
$ws.Cells.Item($row, 2).Value = "Rodillo de recogida de papel RL1-1442-000 para HP"
$ws.Cells.Item($row, 3).Value = "P1005 P1006 P1007 P1008 P1009 P1106 P1108 P1106 P1102 CP1025 M102 M175 M176 M1132 M1136 M1212 M1213 M1214 M1216 M1217"
$ws.Cells.Item($row, 4).Value = 6500
$ws.Cells.Item($row, 5).Value = 50000
$ws.Cells.Item($row, 6).Value = 35
$ws.Cells.Item($row, 7).Value = 31
$ws.Cells.Item($row, 8).Formula = "=(E84-D84)*G84"
$ws.Cells.Item($row, 9).Formula = "=D84*F84"
$ws.Cells.Item($row, 10).Value = 227500
